{"js": "// Locate the \"Exposures\" section paragraph that currently reads:\n//   \"The exposure variable is the audit filters.\"\n// and extend it, then add the bullet list of audit filters below it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"The exposure variable is the audit filters.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph: \" + targetText);\n}\n\n// Append the lead-in sentence to the existing paragraph text.\ntarget.insertText(\" The following audit filters are of interest:\", Word.InsertLocation.end);\n\nconst bulletLines = [\n  \"\u2022 Systolic blood pressure less than 90\",\n  \"\u2022 Glasgow coma scale less than 9 and not intubated\",\n  \"\u2022 Injury severity score greater than 15 but not admitted to the intensive care unit\",\n  \"\u2022 Time to acute intervention more than 60 minutes from arrival to hospital\",\n  \"\u2022 Time to computed tomography more than 30 minutes from arrival to hospital\",\n  \"\u2022 No anticoagulant therapy within 72 hours after traumatic brain injury\",\n  \"\u2022 The presence of cardio-pulmonary resuscitation with thoracotomy\",\n  \"\u2022 The presence of a liver or spleen injury\",\n  \"\u2022 Massive transfusion, defined as 10 or more units of packed red blood cells within 24 hours.\"\n];\n\n// Insert the bullet paragraphs right after the target paragraph, preserving order.\nlet anchor = target;\nfor (const line of bulletLines) {\n  const newPara = anchor.insertParagraph(line, Word.InsertLocation.after);\n  newPara.style = \"Body Text\";\n  anchor = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Locate the \"Exposures\" section paragraph that currently reads:\n#   \"The exposure variable is the audit filters.\"\n# and extend it, then add the bullet list of audit filters below it.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$targetText = \"The exposure variable is the audit filters.\"\n$targetIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  if ($paras.Item($i).Range.Text.Trim() -eq $targetText) {\n    $targetIndex = $i\n    break\n  }\n}\n\nif ($targetIndex -eq -1) {\n  throw \"Could not find target paragraph: $targetText\"\n}\n\n$targetPara = $paras.Item($targetIndex)\n$targetPara.Range.InsertAfter(\" The following audit filters are of interest:\")\n\n$bulletLines = @(\n  \"\u2022 Systolic blood pressure less than 90\",\n  \"\u2022 Glasgow coma scale less than 9 and not intubated\",\n  \"\u2022 Injury severity score greater than 15 but not admitted to the intensive care unit\",\n  \"\u2022 Time to acute intervention more than 60 minutes from arrival to hospital\",\n  \"\u2022 Time to computed tomography more than 30 minutes from arrival to hospital\",\n  \"\u2022 No anticoagulant therapy within 72 hours after traumatic brain injury\",\n  \"\u2022 The presence of cardio-pulmonary resuscitation with thoracotomy\",\n  \"\u2022 The presence of a liver or spleen injury\",\n  \"\u2022 Massive transfusion, defined as 10 or more units of packed red blood cells within 24 hours.\"\n)\n\n$curIndex = $targetIndex\nforeach ($line in $bulletLines) {\n  $curPara = $paras.Item($curIndex)\n  $curPara.Range.InsertParagraphAfter()\n  $curIndex = $curIndex + 1\n  $newPara = $paras.Item($curIndex)\n  $newPara.Range.InsertAfter($line)\n  $newPara.Range.Style = \"Body Text\"\n}\n"}
